$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK_EQ5D")

# Correct spelling of the regressor name in cell A8: "Che_pcs_cb" -> "Dhe_pcs_cb"
$ws.Range("A8").Value = "Dhe_pcs_cb"

# Reflect the resulting cell selection seen in the saved file
$ws.Range("A8").Select()
